# Restore C10 on the "Rules" sheet back to its numeric value of 1
# (was 18, per the target revision).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
